$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.919.09"
$ws.Range("E2").Value = "  -1.70%  "

$ws.Range("D3").Value = "1.832.63"
$ws.Range("E3").Value = "  -1.92%  "

$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "'245.91"
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("D6").Value = "'0.6901"
$ws.Range("E6").Value = "  -2.41%  "

$ws.Range("D7").Value = "'0.9995"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").Value = "'0.07651"
$ws.Range("E8").Value = "  -2.96%  "

$ws.Range("D9").Value = "'0.3054"
$ws.Range("E9").Value = "  -2.52%  "

$ws.Range("D10").Value = "'23.52"
$ws.Range("E10").Value = "  -4.23%  "

$ws.Range("D11").Value = "'0.07822"
$ws.Range("E11").Value = "  -1.52%  "

$ws.Range("D12").Value = "1.834.09"
$ws.Range("E12").Value = "  -2.52%  "

$ws.Range("E13").Value = "  -2.65%  "

$ws.Range("D14").Value = "'90.49"
$ws.Range("E14").Value = "  -3.13%  "

$ws.Range("D15").Value = "'0.6785"
$ws.Range("E15").Value = "  -3.23%  "

$ws.Range("D16").Value = "'6.432"
$ws.Range("E16").Value = "  -1.17%  "

$ws.Range("D17").Value = "'0.000008293"
$ws.Range("E17").Value = "  -0.69%  "

$ws.Range("D18").Value = "28.916.06"
$ws.Range("E18").Value = "  -1.90%  "

$ws.Range("D19").Value = "'242.75"
$ws.Range("E19").Value = "  -3.77%  "

$ws.Range("D20").Value = "2.081.51"
$ws.Range("E20").Value = "  -2.63%  "

$ws.Range("E21").Value = "  -3.39%  "

$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").Value = "'7.448"
$ws.Range("E23").Value = "  -2.46%  "

$ws.Range("D24").Value = "'0.9996"
$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").Value = "'0.1468"
$ws.Range("E25").Value = "  -5.54%  "

$ws.Range("D26").Value = "'161.35"
$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("D27").Value = "'8.797"
$ws.Range("E27").Value = "  -2.25%  "

$ws.Range("D28").Value = "'18.17"
$ws.Range("E28").Value = "  -2.90%  "

$ws.Range("D29").Value = "'1.559"
$ws.Range("E29").Value = "  +3.87%  "

$ws.Range("D30").Value = "'4.216"
$ws.Range("E30").Value = "  -2.64%  "

$ws.Range("D31").Value = "'4.141"
$ws.Range("E31").Value = "  -2.60%  "

$ws.Range("D32").Value = "'1.176"
$ws.Range("E32").Value = "  -2.25%  "

$ws.Range("D33").Value = "'0.05116"
$ws.Range("E33").Value = "  -3.76%  "

$ws.Range("D34").Value = "'0.7559"
$ws.Range("E34").Value = "  +0.96%  "

$ws.Range("D35").Value = "'1.839"
$ws.Range("E35").Value = "  -2.92%  "

$ws.Range("D36").Value = "'1.145"
$ws.Range("E36").Value = "  -2.33%  "

$ws.Range("E37").Value = "  -1.38%  "

$ws.Range("D38").Value = "'0.01841"
$ws.Range("E38").Value = "  -2.58%  "

$ws.Range("D39").Value = "1.231.90"
$ws.Range("E39").Value = "  -3.34%  "

$ws.Range("D40").Value = "'2.688"
$ws.Range("E40").Value = "  -2.35%  "

$ws.Range("D41").Value = "'0.9259"
$ws.Range("E41").Value = "  +3.69%  "

$ws.Range("D42").Value = "'109.03"
$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("D43").Value = "'0.9993"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").Value = "'5.700"
$ws.Range("E44").Value = "  -6.02%  "

$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.980.85"
$ws.Range("E45").Value = "  -2.29%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'9.531"
$ws.Range("E46").Value = "  -0.25%  "

$ws.Range("D47").Value = "'0.5168"
$ws.Range("E47").Value = "  -0.25%  "

$ws.Range("E48").Value = "  -4.52%  "

$ws.Range("D49").Value = "'64.06"
$ws.Range("E49").Value = "  -10.22%  "

$ws.Range("D50").Value = "'1.740"
$ws.Range("E50").Value = "  -3.18%  "

$ws.Range("D51").Value = "'0.4196"
$ws.Range("E51").Value = "  -2.54%  "
